$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.335.17"
$ws.Range("E2").Value = "'  +2.98%  "

# Row 3
$ws.Range("D3").Value = "'3.411.36"
$ws.Range("E3").Value = "'  +2.79%  "

# Row 4
$ws.Range("E4").Value = "'  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'562.34"
$ws.Range("E5").Value = "'  +2.00%  "

# Row 6
$ws.Range("D6").Value = "'175.28"
$ws.Range("E6").Value = "'  +2.70%  "

# Row 7
$ws.Range("E7").Value = "'  +3.00%  "

# Row 8
$ws.Range("D8").Value = "'3.406.90"
$ws.Range("E8").Value = "'  +2.92%  "

# Row 9
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "'  -0.19%  "

# Row 10
$ws.Range("E10").Value = "'  +12.60%  "

# Row 11
$ws.Range("D11").Value = "'0.634"
$ws.Range("E11").Value = "'  +3.45%  "

# Row 12
$ws.Range("D12").Value = "'54.80"
$ws.Range("E12").Value = "'  +3.07%  "

# Row 13
$ws.Range("E13").Value = "'  +5.68%  "

# Row 14
$ws.Range("D14").Value = "'9.17"
$ws.Range("E14").Value = "'  +3.29%  "

# Row 15
$ws.Range("D15").Value = "'3.945.57"
$ws.Range("E15").Value = "'  +2.36%  "

# Row 16
$ws.Range("D16").Value = "'18.42"
$ws.Range("E16").Value = "'  +4.16%  "

# Row 17
$ws.Range("D17").Value = "'3.403.62"
$ws.Range("E17").Value = "'  +2.49%  "

# Row 18
$ws.Range("E18").Value = "'  +1.48%  "

# Row 19
$ws.Range("D19").Value = "'11.93"
$ws.Range("E19").Value = "'  +2.58%  "

# Row 20
$ws.Range("D20").Value = "'65.205.97"
$ws.Range("E20").Value = "'  +2.90%  "

# Row 21
$ws.Range("D21").Value = "'0.997"
$ws.Range("E21").Value = "'  +3.26%  "

# Row 22
$ws.Range("D22").Value = "'473.35"
$ws.Range("E22").Value = "'  +17.12%  "

# Row 23
$ws.Range("D23").Value = "'4.97"
$ws.Range("E23").Value = "'  +16.83%  "

# Row 24
$ws.Range("E24").Value = "'  +2.55%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'86.88"
$ws.Range("E25").Value = "'  +5.25%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'13.80"
$ws.Range("E26").Value = "'  +5.26%  "

# Row 27
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = "'  +3.52%  "

# Row 28
$ws.Range("D28").Value = "'2.90"
$ws.Range("E28").Value = "'  +6.67%  "

# Row 29
$ws.Range("D29").Value = "'8.88"
$ws.Range("E29").Value = "'  +3.85%  "

# Row 30
$ws.Range("D30").Value = "'30.78"
$ws.Range("E30").Value = "'  +6.01%  "

# Row 31
$ws.Range("D31").Value = "'6.72"
$ws.Range("E31").Value = "'  +3.87%  "

# Row 32
$ws.Range("D32").Value = "'11.58"
$ws.Range("E32").Value = "'  +2.81%  "

# Row 33
$ws.Range("D33").Value = "'585.15"
$ws.Range("E33").Value = "'  +1.90%  "

# Row 34
$ws.Range("E34").Value = "'  +3.56%  "

# Row 35
$ws.Range("D35").Value = "'60.46"
$ws.Range("E35").Value = "'  +5.33%  "

# Row 36
$ws.Range("E36").Value = "'  -0.13%  "

# Row 37
$ws.Range("E37").Value = "'  -3.66%  "

# Row 38
$ws.Range("D38").Value = "'3.51"
$ws.Range("E38").Value = "'  +3.53%  "

# Row 39
$ws.Range("D39").Value = "'36.04"
$ws.Range("E39").Value = "'  +3.28%  "

# Row 40
$ws.Range("D40").Value = "'0.0₃0752"
$ws.Range("E40").Value = "'  +2.81%  "

# Row 41
$ws.Range("D41").Value = "'0.375"
$ws.Range("E41").Value = "'  +2.91%  "

# Row 42
$ws.Range("D42").Value = "'3.113.66"
$ws.Range("E42").Value = "'  -0.80%  "

# Row 43
$ws.Range("E43").Value = "'  +0.16%  "

# Row 44
$ws.Range("D44").Value = "'2.87"
$ws.Range("E44").Value = "'  +2.72%  "

# Row 45
$ws.Range("D45").Value = "'2.53"
$ws.Range("E45").Value = "'  +3.95%  "

# Row 46
$ws.Range("E46").Value = "'  +3.69%  "

# Row 47
$ws.Range("D47").Value = "'3.21"
$ws.Range("E47").Value = "'  +1.57%  "

# Row 48
$ws.Range("E48").Value = "'  +5.53%  "

# Row 49
$ws.Range("D49").Value = "'2.58"
$ws.Range("E49").Value = "'  -0.42%  "

# Row 50
$ws.Range("D50").Value = "'137.08"
$ws.Range("E50").Value = "'  +3.32%  "

# Row 51
$ws.Range("D51").Value = "'8.38"
$ws.Range("E51").Value = "'  +4.85%  "
